$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated "Max Decel" curve values for B2:B32 (Speed 0..30 m/s)
# The curve now decelerates faster at first and plateaus at -1.03162005421553
$values = @(
    -2.4230514575044602,
    -2.1869262666181899,
    -1.9524242114449399,
    -1.71959152131576,
    -1.48847749822579,
    -1.2591348467896999,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553,
    -1.03162005421553
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Range("B$row").Value = $values[$i]
}

# The chart now only plots the first 24 points (Sheet1!$B$2:$B$25)
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Values = "=Sheet1!`$B`$2:`$B`$25"
